$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # New headers in columns G and H
    $ws.Range("G1").Value = "p_ttes_cop"
    $ws.Range("H1").Value = "p_ttes_c_charge_discharge"

    # Replace the formula in E2 with a plain literal value
    $ws.Range("E2").Value = 2500

    # New data values in columns G and H
    $ws.Range("G2").Value = 5.3
    $ws.Range("H2").Value = 5
}
